# Update countries & provincias Spain
# Applies the 9-Sep-2020 14:28 data refresh to the "Pais" sheet:
#   - Ghana overtakes Armenia in total cases (rows 61/62 swap country order)
#   - Gambia overtakes Siria in total cases (rows 130/131 swap country order)
#   - Liechtenstein overtakes Curazao in total cases (rows 194/195 swap country order)
#   - Refreshed case/recovered/active/critical/death counters for several countries
#   - "Datos actualizados" timestamp bumped from 13:11 to 14:28

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 14:28"

# --- Estados Unidos (row 4): unchanged country, refreshed counters ---
$ws.Cells.Item(4,2).Value = 6514433
$ws.Cells.Item(4,3).Value = 202
$ws.Cells.Item(4,4).Value = 3797174
$ws.Cells.Item(4,5).Value = 2523222

# --- Kuwait (row 38): unchanged country, refreshed counters ---
$ws.Cells.Item(38,2).Value = 92082
$ws.Cells.Item(38,3).Value = 838
$ws.Cells.Item(38,4).Value = 82222
$ws.Cells.Item(38,5).Value = 9308
$ws.Cells.Item(38,7).Value = 4
$ws.Cells.Item(38,8).Value = 552

# --- Suiza (row 60): unchanged country, refreshed counters ---
$ws.Cells.Item(60,4).Value = 38100
$ws.Cells.Item(60,5).Value = 5188

# --- Ghana now outranks Armenia: row 61 becomes Ghana, row 62 becomes Armenia ---
$ws.Cells.Item(61,1).Value = "Ghana"
$ws.Cells.Item(61,2).Value = 45188
$ws.Cells.Item(61,3).Value = 176
$ws.Cells.Item(61,4).Value = 44042
$ws.Cells.Item(61,5).Value = 863
$ws.Cells.Item(61,7).Value = 0
$ws.Cells.Item(61,8).Value = 283

$ws.Cells.Item(62,1).Value = "Armenia"
$ws.Cells.Item(62,2).Value = 45152
$ws.Cells.Item(62,3).Value = 199
$ws.Cells.Item(62,4).Value = 41023
$ws.Cells.Item(62,5).Value = 3224
$ws.Cells.Item(62,7).Value = 2
$ws.Cells.Item(62,8).Value = 905

# --- Kenia (row 72): unchanged country, refreshed counters ---
$ws.Cells.Item(72,4).Value = 20164
$ws.Cells.Item(72,5).Value = 9272

# --- Serbia (row 73): unchanged country, refreshed counters ---
$ws.Cells.Item(73,2).Value = 27919
$ws.Cells.Item(73,3).Value = 556
$ws.Cells.Item(73,4).Value = 18466
$ws.Cells.Item(73,5).Value = 9261
$ws.Cells.Item(73,7).Value = 2
$ws.Cells.Item(73,8).Value = 192

# --- Austria (row 74): unchanged country, refreshed counters ---
$ws.Cells.Item(74,2).Value = 26602
$ws.Cells.Item(74,3).Value = 91
$ws.Cells.Item(74,4).Value = 16786
$ws.Cells.Item(74,5).Value = 9046

# --- Madagascar (row 85): unchanged country, refreshed counters ---
$ws.Cells.Item(85,2).Value = 15520
$ws.Cells.Item(85,3).Value = 85
$ws.Cells.Item(85,4).Value = 14243
$ws.Cells.Item(85,5).Value = 1071

# --- Row 103: unchanged country, refreshed counters ---
$ws.Cells.Item(103,5).Value = 743
$ws.Cells.Item(103,7).Value = 1
$ws.Cells.Item(103,8).Value = 337

# --- Gambia now outranks Siria: row 130 becomes Gambia, row 131 becomes Siria ---
$ws.Cells.Item(130,1).Value = "Gambia"
$ws.Cells.Item(130,2).Value = 3293
$ws.Cells.Item(130,3).Value = 18
$ws.Cells.Item(130,4).Value = 1460
$ws.Cells.Item(130,5).Value = 1734
$ws.Cells.Item(130,8).Value = 99

$ws.Cells.Item(131,1).Value = "Siria"
$ws.Cells.Item(131,2).Value = 3289
$ws.Cells.Item(131,4).Value = 760
$ws.Cells.Item(131,5).Value = 2389
$ws.Cells.Item(131,8).Value = 140

# --- Liechtenstein now outranks Curazao: row 194 becomes Liechtenstein, row 195 becomes Curazao ---
$ws.Cells.Item(194,1).Value = "Liechtenstein"
$ws.Cells.Item(194,2).Value = 108
$ws.Cells.Item(194,3).Value = 1
$ws.Cells.Item(194,4).Value = 105
$ws.Cells.Item(194,5).Value = 2

$ws.Cells.Item(195,1).Value = "Curazao"
$ws.Cells.Item(195,4).Value = 48
$ws.Cells.Item(195,5).Value = 58
